$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns C and D keep their original text (string) storage format
# by forcing a text number format before assigning the new values, so that
# numeric-looking strings are not silently re-typed as numbers by Excel.
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "157"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "416386.40"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "452"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1167510.82"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "221"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "483041.00"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "814"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3097688.81"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "34"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "94600.00"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "176"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "542816.18"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "104"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "250800.00"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "98"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "237788.98"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "137"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "626217.26"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "194"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "440089.87"
$ws.Range("C77").NumberFormat = "@"
$ws.Range("C77").Value = "91"
$ws.Range("D77").NumberFormat = "@"
$ws.Range("D77").Value = "240487.00"
$ws.Range("C78").NumberFormat = "@"
$ws.Range("C78").Value = "217"
$ws.Range("D78").NumberFormat = "@"
$ws.Range("D78").Value = "605075.19"
$ws.Range("C80").NumberFormat = "@"
$ws.Range("C80").Value = "500"
$ws.Range("D80").NumberFormat = "@"
$ws.Range("D80").Value = "2198776.03"
$ws.Range("C122").NumberFormat = "@"
$ws.Range("C122").Value = "256"
$ws.Range("D122").NumberFormat = "@"
$ws.Range("D122").Value = "712508.00"
$ws.Range("C123").NumberFormat = "@"
$ws.Range("C123").Value = "128"
$ws.Range("D123").NumberFormat = "@"
$ws.Range("D123").Value = "331012.45"
$ws.Range("C124").NumberFormat = "@"
$ws.Range("C124").Value = "515"
$ws.Range("D124").NumberFormat = "@"
$ws.Range("D124").Value = "2345836.06"
$ws.Range("C132").NumberFormat = "@"
$ws.Range("C132").Value = "89"
$ws.Range("D132").NumberFormat = "@"
$ws.Range("D132").Value = "396663.75"
$ws.Range("C138").NumberFormat = "@"
$ws.Range("C138").Value = "578"
$ws.Range("D138").NumberFormat = "@"
$ws.Range("D138").Value = "1457046.00"
$ws.Range("C140").NumberFormat = "@"
$ws.Range("C140").Value = "2724"
$ws.Range("D140").NumberFormat = "@"
$ws.Range("D140").Value = "6894397.55"
$ws.Range("C141").NumberFormat = "@"
$ws.Range("C141").Value = "2700"
$ws.Range("D141").NumberFormat = "@"
$ws.Range("D141").Value = "11816211.82"
$ws.Range("C145").NumberFormat = "@"
$ws.Range("C145").Value = "1066"
$ws.Range("D145").NumberFormat = "@"
$ws.Range("D145").Value = "2805349.25"
$ws.Range("C149").NumberFormat = "@"
$ws.Range("C149").Value = "436"
$ws.Range("D149").NumberFormat = "@"
$ws.Range("D149").Value = "1433905.46"
$ws.Range("C150").NumberFormat = "@"
$ws.Range("C150").Value = "865"
$ws.Range("D150").NumberFormat = "@"
$ws.Range("D150").Value = "2106195.82"
$ws.Range("C157").NumberFormat = "@"
$ws.Range("C157").Value = "8"
$ws.Range("D157").NumberFormat = "@"
$ws.Range("D157").Value = "21900.00"
$ws.Range("C201").NumberFormat = "@"
$ws.Range("C201").Value = "669"
$ws.Range("D201").NumberFormat = "@"
$ws.Range("D201").Value = "2567782.58"
$ws.Range("C202").NumberFormat = "@"
$ws.Range("C202").Value = "27"
$ws.Range("D202").NumberFormat = "@"
$ws.Range("D202").Value = "100238.00"
$ws.Range("C207").NumberFormat = "@"
$ws.Range("C207").Value = "81"
$ws.Range("D207").NumberFormat = "@"
$ws.Range("D207").Value = "192720.00"
$ws.Range("C210").NumberFormat = "@"
$ws.Range("C210").Value = "140"
$ws.Range("D210").NumberFormat = "@"
$ws.Range("D210").Value = "323306.36"
$ws.Range("C226").NumberFormat = "@"
$ws.Range("C226").Value = "19"
$ws.Range("D226").NumberFormat = "@"
$ws.Range("D226").Value = "60578.00"
$ws.Range("C228").NumberFormat = "@"
$ws.Range("C228").Value = "59"
$ws.Range("D228").NumberFormat = "@"
$ws.Range("D228").Value = "155500.00"
$ws.Range("C229").NumberFormat = "@"
$ws.Range("C229").Value = "166"
$ws.Range("D229").NumberFormat = "@"
$ws.Range("D229").Value = "461905.00"
$ws.Range("C230").NumberFormat = "@"
$ws.Range("C230").Value = "21"
$ws.Range("D230").NumberFormat = "@"
$ws.Range("D230").Value = "57403.00"
$ws.Range("C231").NumberFormat = "@"
$ws.Range("C231").Value = "345"
$ws.Range("D231").NumberFormat = "@"
$ws.Range("D231").Value = "1213037.07"
$ws.Range("C232").NumberFormat = "@"
$ws.Range("C232").Value = "7"
$ws.Range("D232").NumberFormat = "@"
$ws.Range("D232").Value = "31972.00"
$ws.Range("C234").NumberFormat = "@"
$ws.Range("C234").Value = "19"
$ws.Range("D234").NumberFormat = "@"
$ws.Range("D234").Value = "40612.00"
$ws.Range("C235").NumberFormat = "@"
$ws.Range("C235").Value = "83"
$ws.Range("D235").NumberFormat = "@"
$ws.Range("D235").Value = "250687.09"
$ws.Range("C236").NumberFormat = "@"
$ws.Range("C236").Value = "29"
$ws.Range("D236").NumberFormat = "@"
$ws.Range("D236").Value = "90323.00"
$ws.Range("C237").NumberFormat = "@"
$ws.Range("C237").Value = "31"
$ws.Range("D237").NumberFormat = "@"
$ws.Range("D237").Value = "89100.00"
$ws.Range("C239").NumberFormat = "@"
$ws.Range("C239").Value = "66"
$ws.Range("D239").NumberFormat = "@"
$ws.Range("D239").Value = "291949.36"
$ws.Range("C240").NumberFormat = "@"
$ws.Range("C240").Value = "71"
$ws.Range("D240").NumberFormat = "@"
$ws.Range("D240").Value = "149500.00"
$ws.Range("C243").NumberFormat = "@"
$ws.Range("C243").Value = "164"
$ws.Range("D243").NumberFormat = "@"
$ws.Range("D243").Value = "419700.00"
$ws.Range("C244").NumberFormat = "@"
$ws.Range("C244").Value = "541"
$ws.Range("D244").NumberFormat = "@"
$ws.Range("D244").Value = "1410575.83"
$ws.Range("C245").NumberFormat = "@"
$ws.Range("C245").Value = "117"
$ws.Range("D245").NumberFormat = "@"
$ws.Range("D245").Value = "342727.11"
$ws.Range("C246").NumberFormat = "@"
$ws.Range("C246").Value = "1016"
$ws.Range("D246").NumberFormat = "@"
$ws.Range("D246").Value = "3766966.46"
$ws.Range("C247").NumberFormat = "@"
$ws.Range("C247").Value = "44"
$ws.Range("D247").NumberFormat = "@"
$ws.Range("D247").Value = "116077.85"
$ws.Range("C248").NumberFormat = "@"
$ws.Range("C248").Value = "25"
$ws.Range("D248").NumberFormat = "@"
$ws.Range("D248").Value = "62000.00"
$ws.Range("C249").NumberFormat = "@"
$ws.Range("C249").Value = "83"
$ws.Range("D249").NumberFormat = "@"
$ws.Range("D249").Value = "188500.00"
$ws.Range("C250").NumberFormat = "@"
$ws.Range("C250").Value = "187"
$ws.Range("D250").NumberFormat = "@"
$ws.Range("D250").Value = "587139.19"
$ws.Range("C251").NumberFormat = "@"
$ws.Range("C251").Value = "131"
$ws.Range("D251").NumberFormat = "@"
$ws.Range("D251").Value = "414193.00"
$ws.Range("C252").NumberFormat = "@"
$ws.Range("C252").Value = "103"
$ws.Range("D252").NumberFormat = "@"
$ws.Range("D252").Value = "273972.92"
$ws.Range("C254").NumberFormat = "@"
$ws.Range("C254").Value = "137"
$ws.Range("D254").NumberFormat = "@"
$ws.Range("D254").Value = "502835.82"
$ws.Range("C255").NumberFormat = "@"
$ws.Range("C255").Value = "222"
$ws.Range("D255").NumberFormat = "@"
$ws.Range("D255").Value = "495863.00"
